$wb = $excel.ActiveWorkbook
$ws4 = $wb.Worksheets.Item(4)
$ws4.Columns.Item(9).ColumnWidth = 41
Write-Host "ok"
